$wb = $excel.ActiveWorkbook

# --- DHW sheet ---
$dhw = $wb.Worksheets.Item("DHW")
$dhw.Range("C3").Value = "BOILER"
$dhw.Range("C4").Value = "BOILER"
$dhw.Range("C5").Value = "BOILER"
$dhw.Range("C7").Value = "BOILER"
$dhw.Range("C9").Value = "SC"
$dhw.Range("C16").Value = "SC"
$dhw.Range("C17").Value = "SC"
$dhw.Columns("E").Select()

# --- HEATING sheet (active sheet) ---
$heating = $wb.Worksheets.Item("HEATING")
$heating.Activate()
$heating.Range("C3").Value = "BOILER"
$heating.Range("C4").Value = "BOILER"
$heating.Range("C5").Value = "BOILER"
$heating.Range("C7").Value = "BOILER"
$heating.Range("C9").Value = "SC"
$heating.Range("C16").Value = "SC"
$heating.Range("C17").Value = "SC"

# Match formatting of E7 to the equivalent DHW cell (picks up the "Normal 2 2" cell style)
$dhw.Range("E7").Copy()
$heating.Range("E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$heating.Columns("E").Select()

# --- COOLING sheet ---
$cooling = $wb.Worksheets.Item("COOLING")
$cooling.Range("C2").Select()

# --- ELECTRICITY sheet ---
$electricity = $wb.Worksheets.Item("ELECTRICITY")
$electricity.Range("C8").Formula = "=HEATING!D5/0.4"

$heating.Activate()
